$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column B label "age" -> "Age"
$ws.Range("B1").Value = "Age"

# Row 2 (A2:E2) is unchanged in the source data.

# New data table (rows 3-52) reflecting the refreshed "all_files_01_2024" data frame.
# Each row is written as RegistryCode, Age, WSO, incongruent, congruent.
$rows = @(
    @(3, "sub 3058", 61, 1422.428571428571, 0.7158469945355191, 0.6833333333333333),
    @(4, "sub 2917", 59, 450.8571428571428, 0.9526315789473684, 0.9105263157894737),
    @(5, "sub 3104", 50, 405.5714285714286, 0.7671957671957672, 0.7329842931937173),
    @(6, "sub 3889", 62, 202.7142857142857, 0.8128342245989305, 0.7936507936507936),
    @(7, "sub 3642", 60, 275, 0.8, 0.75),
    @(8, "sub 3035", 67, 843.2857142857143, 0.7248677248677249, 0.7291666666666666),
    @(9, "sub 4182", 53, 146.2857142857143, 0.9, 0.8617021276595744),
    @(10, "sub 4137", 70, 162.2857142857143, 0.78125, 0.6649214659685864),
    @(11, "sub 3583", 51, 481.5714285714286, 0.4450261780104712, 0.5631578947368421),
    @(12, "sub 4281", 46, 128.5714285714286, 0.7368421052631579, 0.7),
    @(13, "sub 4191", 55, 131.4285714285714, 0.7015706806282722, 0.680628272251309),
    @(14, "sub 3201", 51, 436.4285714285714, 0.6858638743455497, 0.6157894736842106),
    @(15, "sub 4275", 55, 121.5714285714286, 0.6042780748663101, 0.5873015873015873),
    @(16, "sub 4170", 58, 308, 0.7407407407407407, 0.7591623036649214),
    @(17, "sub 3676", 63, 280.8571428571428, 0.7947368421052632, 0.7724867724867724),
    @(18, "sub 4431", 71, 96.71428571428571, 0.5925925925925926, 0.6894736842105263),
    @(19, "sub 4294", 67, 107.5714285714286, 0.6276595744680851, 0.6417112299465241),
    @(20, "sub 3650", 62, 331.5714285714286, 0.7058823529411765, 0.6473684210526316),
    @(21, "sub 4140", 44, 161.7142857142857, 0.7010869565217391, 0.7021276595744681),
    @(22, "sub 3912", 55, 229.8571428571429, 0.6736842105263158, 0.6105263157894737),
    @(23, "sub 4051", 59, 196.7142857142857, 0.6412213740458015, 0.810126582278481),
    @(24, "sub 3286", 76, 390.1428571428572, 0.875, 0.8085106382978723),
    @(25, "sub 4208", 51, 137.2857142857143, 0.7329842931937173, 0.7604166666666666),
    @(26, "sub 1536", 69, 1061, 0.6878306878306878, 0.6387434554973822),
    @(27, "sub 4439", 56, 79.14285714285714, 0.7248677248677249, 0.6666666666666666),
    @(28, "sub 4466", 65, 95.85714285714286, 0.774869109947644, 0.6666666666666666),
    @(29, "sub 4504", 55, 78.14285714285714, 0.6935483870967742, 0.612565445026178),
    @(30, "sub 3887", 68, 310.4285714285714, 0.9214659685863874, 0.7736842105263158),
    @(31, "sub 4391", 61, 101, 0.723404255319149, 0.7368421052631579),
    @(32, "sub 4546", 40, 123.5714285714286, 0.9791666666666666, 0.9635416666666666),
    @(33, "sub 4227", 38, 186.7142857142857, 0, 1),
    @(34, "sub 4532", 63, 89.42857142857143, 0.7291666666666666, 0.8219895287958116),
    @(35, "sub 3299", 60, 509.2857142857143, 0.8723404255319149, 0.859375),
    @(36, "sub 4720", 59, 43.28571428571428, 0.7142857142857143, 0.6666666666666666),
    @(37, "sub 4607", 63, 60.42857142857143, 0.7513227513227513, 0.4947368421052631),
    @(38, "sub 4728", 43, 34.42857142857143, 0.7853403141361257, 0.7801047120418848),
    @(39, "sub 4769", 54, 31.71428571428572, 0.632768361581921, 0.6171428571428571),
    @(40, "sub 4777", 51, 28.71428571428572, 0.8052631578947368, 0.7473684210526316),
    @(41, "sub 3396", 54, 337.2857142857143, 0.5531914893617021, 0.5833333333333334),
    @(42, "sub 4467", 60, 31.28571428571428, 0.6170212765957447, 0.4516129032258064),
    @(43, "sub 3392", 53, 331.1428571428572, 0.6387434554973822, 0.6230366492146597),
    @(44, "sub 3960", 79, 239.7142857142857, 0.3397435897435898, 0.6329113924050633),
    @(45, "sub 3154", 30, 443.1428571428572, 0.8010471204188482, 0.746031746031746),
    @(46, "sub 3994", 55, 179.8571428571429, 0.4867724867724867, 0.4712041884816754),
    @(47, "sub 4086", 59, 141.7142857142857, 0.921875, 0.9157894736842105),
    @(48, "sub 4002", 69, 185.2857142857143, 0.6595744680851063, 0.7864583333333334),
    @(49, "sub 2522", 40, 555.8571428571429, 0.8624338624338624, 0.765625),
    @(50, "sub 3364", 34, 334.4285714285714, 0.9267015706806283, 0.6577540106951871),
    @(51, "sub 2884", 55, 516.8571428571429, 0.6063829787234043, 0.5661375661375662),
    @(52, "sub 2998", 59, 469.2857142857143, 0.8926553672316384, 0.8882978723404256)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Output "Updated header and rows 3-52 with refreshed data."